# Apply the RPAR_holdings.xlsx update:
#  1. Bump the "as of" date in the confidential disclaimer note from
#     2021-05-24 to 2021-05-25.
#  2. Refresh the Weight (D) and Percent Change (E) values for rows 2-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected, so the cells have to be unlocked for editing
# first; re-protect at the end to leave the workbook in the same
# (protected) state it was delivered in.
$ws.Unprotect("#ti!8*!")

# 1) Update the disclaimer text (A18) with the new "as of" date.
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-25 for illustrative purposes only and are subject to change."

# 2) Update Weight / Percent Change figures for each holding row.
$ws.Range("D2").Value = 0.05691396264841002
$ws.Range("E2").Value = -0.002763194252555889

$ws.Range("D3").Value = 0.02392973429535421
$ws.Range("E3").Value = 0

$ws.Range("D4").Value = 0.03099680032606139
$ws.Range("E4").Value = 0.01046423135464236

$ws.Range("D5").Value = 0.03280599544064484
$ws.Range("E5").Value = -0.01821631878557872

$ws.Range("D6").Value = 0.03677233300984934
$ws.Range("E6").Value = -0.01892440928044481

$ws.Range("D7").Value = 0.01882912771371076
$ws.Range("E7").Value = -0.004384896467722221

$ws.Range("D8").Value = 0.004446649057429659
$ws.Range("E8").Value = -0.002652519893899252

$ws.Range("D9").Value = 0.006884246241873857
$ws.Range("E9").Value = 0.001142204454597451

$ws.Range("D10").Value = 0.07359970853676676
$ws.Range("E10").Value = 0.008547008547008739

$ws.Range("D11").Value = 0.07363902462038682
$ws.Range("E11").Value = 0.009610250934329878

$ws.Range("D12").Value = 0.1448719049232255
$ws.Range("E12").Value = 0.009263279779997191

$ws.Range("D13").Value = 0.382042247752922
$ws.Range("E13").Value = 0.003673897830650752

$ws.Range("D14").Value = 0.1142682654333648
$ws.Range("E14").Value = 0.009995183044315903

$ws.Range("E15").Value = 0.004011551065367547

# Re-protect the sheet with the same protection options it had before.
$ws.Protect("#ti!8*!", $true, $true, $true)
